$wb = $excel.ActiveWorkbook

# Sheets "展览" (sheet1) and "全部类型" (sheet4) contain the same data set
# and both need the identical updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 167
    $ws.Range("F3").Value = 426
    $ws.Range("F4").Value = 12361
    $ws.Range("F6").Value = 141
    $ws.Range("F10").Value = 192
    $ws.Range("F11").Value = 449
    $ws.Range("F12").Value = 59
    $ws.Range("F15").Value = 42

    $ws.Range("C16").Value = "合肥·运动番only·群青日和"
    $ws.Range("F16").Value = 369
    $ws.Range("I16").Value = "//i2.hdslb.com/bfs/openplatform/202404/Jzeq47lD1714026878824.jpeg"

    $ws.Range("F17").Value = 3447
    $ws.Range("F19").Value = 937
    $ws.Range("F22").Value = 38
}
